# Cambiar a tipo float datos con decimales
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the fill color used by the header style (fill index 3: FFFAEBD7 -> FFEFBC87)
#    This style (s="2") is applied to the header row cells (A1:H1 ...).
$headerRange = $ws.Range("A1:H1")
$headerRange.Interior.Color = 8895727   # RGB(135,188,239) == 0xEFBC87 little-endian OLE color

# 2. Column width changes (BU=8, BV=10 OOXML "width" units).
#    The COM ColumnWidth unit differs from the stored OOXML width by a
#    pixel-rounding step (MDW=6 here): stored = round(ColumnWidth*6)/6 + 5/6.
#    Solve for the ColumnWidth that lands exactly on the target stored width.
$ws.Columns.Item(73).ColumnWidth = 7.166666666666667   # -> stored width 8
$ws.Columns.Item(74).ColumnWidth = 9.166666666666666   # -> stored width 10

# 3. Row 3 cell value changes (convert integers to floats with decimals)
$ws.Range("BL3").Value2 = 54000.75
$ws.Range("BP3").Value2 = 875.05
$ws.Range("BS3").Value2 = 1345.26
$ws.Range("BU3").Value2 = 123.24
$ws.Range("BV3").Value2 = 76547.76700000001
$ws.Range("BW3").Value2 = 1234.45
$ws.Range("CA3").Value2 = 200000.5

# 4. Row 4 cell value changes
$ws.Range("BU4").Value2 = 4645.59
$ws.Range("BV4").Value2 = 98651.25
$ws.Range("BW4").Value2 = 456757.57
$ws.Range("CA4").Value2 = 85.01000000000001
